# Applies the "Updated to include Pathfinder Project" edit to the cover letter.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $find)
    }
    return $ok
}

# 1. Tighten a couple of double-spaces into single spaces.
Replace-Text "I am writing to be considered for a position as a junior software developer.  " "I am writing to be considered for a position as a junior software developer. "
Replace-Text "in computer science.  Since that decision" "in computer science. Since that decision"
Replace-Text "ccessful software developer.  I completed MIT" "ccessful software developer. I completed MIT"

# 2. "...taught in Python and worked through an introductory C++ textbook.  I read..."
#    -> "...taught in Python, and I read..."
Replace-Text "taught in Python and worked through an introductory C++ textbook.  I read" "taught in Python, and I read"

# 3. "...by Goodrich and Tamassia and completed relevant exercises..."
#    -> "...by Goodrich and Tamassia, completing relevant exercises..."
Replace-Text "by Goodrich and Tamassia and completed relevant exercises" "by Goodrich and Tamassia, completing relevant exercises"

# 4. Rework the GitHub / TEAMMATES / Pathfinder paragraph (also absorbs the
#    trailing whitespace that used to sit before the _GoBack bookmark).
Replace-Text "to practice implementing newly learned abstract data types. I have uploaded many of my Java exercises to Github.  Most recently, I contributed to TEAMMATES, open-source software that manages student/instructor correspondence in higher education.   " "to practice implementing newly learned abstract data types. I recently contributed to TEAMMATES, open-source software that manages student/instructor correspondence in higher education, and I am currently developing a Pathfinder Character Builder. My project repositories, including significant commits, are accessible on my Github (https://github.com/carsonshoupe)."

# 5. Rework the "In the coming months" goals paragraph to mention the Pathfinder project.
Replace-Text "In the coming months, my short-term goals are to continue to study foundational texts in computer science literature and make further contributions to TEAMMATES.  I am currently reading sections of " "In the coming weeks, my short-term goals are to complete my Pathfinder Character Builder and continue to study foundational texts in computer science literature.  I designed my character builder to utilize the MVC architectural pattern.  I am near completing the backend of my character builder and next plan to implement the frontend and control.  I am currently reading sections of "

# 6. Word's hidden "_GoBack" bookmark (marks the last edit location) moves from the
#    end of the Github/Pathfinder paragraph into the middle of "I am currently
#    reading sections of ..." (right after "I a"), matching the authored edit.
$rng = $d.Content
$found = $rng.Find.Execute("I am currently reading sections of", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmPoint = $rng.Start + 3
    $bmRange = $d.Range($bmPoint, $bmPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
} else {
    Write-Output "BOOKMARK ANCHOR NOT FOUND"
}

Write-Output "done"
